$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0

$ws.Range("H61").Value = 422.5
$ws.Range("I61").Value = 422.5
$ws.Range("K61").Value = 1267.5
$ws.Range("M61").Value = -1095.5

$ws.Range("H106").Value = 35716260
$ws.Range("I106").Value = 55557030
$ws.Range("J106").Value = 2870.1
$ws.Range("K106").Value = 55557030
$ws.Range("L106").Value = 2870.1
$ws.Range("M106").Value = -55556399
$ws.Range("N106").Value = -4132.1

$ws.Range("H116").Value = 29178042
$ws.Range("J116").Value = 33350134
$ws.Range("L116").Value = 33350134
$ws.Range("N116").Value = -33357018

$ws.Range("H132").Value = 2365222.8
$ws.Range("I132").Value = 1203.5778
$ws.Range("J132").Value = 55555656
$ws.Range("K132").Value = 3610.7334
$ws.Range("L132").Value = 166666968
$ws.Range("M132").Value = -1080.7334
$ws.Range("N132").Value = -166672028

$ws.Range("H137").Value = 78187330
$ws.Range("I137").Value = 1300
$ws.Range("J137").Value = 104249330
$ws.Range("K137").Value = 3900
$ws.Range("L137").Value = 312747990
$ws.Range("M137").Value = -1350
$ws.Range("N137").Value = -312753090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2871917.2
$ws.Range("I61").Value = 1263411.2
$ws.Range("K61").Value = 1263411.2
$ws.Range("M61").Value = -1263199.2

$ws.Range("H74").Value = 91039870
$ws.Range("I74").Value = 92209260
$ws.Range("J74").Value = 88896000
$ws.Range("K74").Value = 92209260
$ws.Range("L74").Value = 88896000
$ws.Range("M74").Value = -92208386
$ws.Range("N74").Value = -88897748

$ws.Range("H77").Value = 91039870
$ws.Range("I77").Value = 92209260
$ws.Range("J77").Value = 88896000
$ws.Range("K77").Value = 461046300
$ws.Range("L77").Value = 444480000
$ws.Range("M77").Value = -461041932
$ws.Range("N77").Value = -444488736

$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

$ws.Range("H132").Value = 11441971
$ws.Range("I132").Value = 11498596
$ws.Range("J132").Value = 11113551
$ws.Range("K132").Value = 34495788
$ws.Range("L132").Value = 33340653
$ws.Range("M132").Value = -34493258
$ws.Range("N132").Value = -33345713

$ws.Range("H136").Value = 2871917.2
$ws.Range("I136").Value = 1263411.2
$ws.Range("K136").Value = 3790233.6
$ws.Range("M136").Value = -3787683.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1600
$ws.Range("I11").Value = 1600
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1600
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1460

$ws.Range("H134").Value = 7050103
$ws.Range("I134").Value = 7576762.5
$ws.Range("J134").Value = 3574151.2
$ws.Range("K134").Value = 22730287.5
$ws.Range("L134").Value = 10722453.6
$ws.Range("M134").Value = -22727752.5
$ws.Range("N134").Value = -10727523.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H31").Value = 3160579.8
$ws.Range("I31").Value = 1985917.8
$ws.Range("J31").Value = 5216238.5
$ws.Range("K31").Value = 1985917.8
$ws.Range("L31").Value = 5216238.5
$ws.Range("M31").Value = -1985622.8
$ws.Range("N31").Value = -5216828.5

$ws.Range("H34").Value = 3160579.8
$ws.Range("I34").Value = 1985917.8
$ws.Range("J34").Value = 5216238.5
$ws.Range("K34").Value = 1985917.8
$ws.Range("L34").Value = 5216238.5
$ws.Range("M34").Value = -1985715.8
$ws.Range("N34").Value = -5216642.5

$ws.Range("H58").Value = 4062590
$ws.Range("I58").Value = 2236859.5
$ws.Range("K58").Value = 2236859.5
$ws.Range("M58").Value = -2236656.5

$ws.Range("H62").Value = 2868.4211
$ws.Range("I62").Value = 2460
$ws.Range("J62").Value = 4400
$ws.Range("K62").Value = 2460
$ws.Range("L62").Value = 4400
$ws.Range("M62").Value = -1836
$ws.Range("N62").Value = -5648

$ws.Range("H65").Value = 2868.4211
$ws.Range("I65").Value = 2460
$ws.Range("J65").Value = 4400
$ws.Range("K65").Value = 12300
$ws.Range("L65").Value = 22000
$ws.Range("M65").Value = -9180
$ws.Range("N65").Value = -28240

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws.Range("H132").Value = 1925252.6
$ws.Range("I132").Value = 2633033.5
$ws.Range("J132").Value = 4133.143
$ws.Range("K132").Value = 7899100.5
$ws.Range("L132").Value = 12399.429
$ws.Range("M132").Value = -7896570.5
$ws.Range("N132").Value = -17459.429

$ws.Range("H134").Value = 1339040.8
$ws.Range("I134").Value = 6259.05
$ws.Range("J134").Value = 4004604.2
$ws.Range("K134").Value = 18777.15
$ws.Range("L134").Value = 12013812.6
$ws.Range("M134").Value = -16242.15
$ws.Range("N134").Value = -12018882.6

$ws.Range("H136").Value = 4062590
$ws.Range("I136").Value = 2236859.5
$ws.Range("K136").Value = 6710578.5
$ws.Range("M136").Value = -6708028.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 166667630
$ws.Range("I120").Value = 166667630
$ws.Range("K120").Value = 500002890
$ws.Range("M120").Value = -499998052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2224995.8
$ws.Range("I70").Value = 1548860.5
$ws.Range("J70").Value = 3036357.8
$ws.Range("K70").Value = 1548860.5
$ws.Range("L70").Value = 3036357.8
$ws.Range("M70").Value = -1548590.5
$ws.Range("N70").Value = -3036897.8

$ws.Range("H73").Value = 2224995.8
$ws.Range("I73").Value = 1548860.5
$ws.Range("J73").Value = 3036357.8
$ws.Range("K73").Value = 1548860.5
$ws.Range("L73").Value = 3036357.8
$ws.Range("M73").Value = -1547924.5
$ws.Range("N73").Value = -3038229.8

$ws.Range("H132").Value = 13022650
$ws.Range("I132").Value = 15477520
$ws.Range("J132").Value = 9094860
$ws.Range("K132").Value = 46432560
$ws.Range("L132").Value = 27284580
$ws.Range("M132").Value = -46430030
$ws.Range("N132").Value = -27289640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1615.3448
$ws.Range("I40").Value = 991.61536
$ws.Range("J40").Value = 2122.125
$ws.Range("K40").Value = 991.61536
$ws.Range("L40").Value = 2122.125
$ws.Range("M40").Value = -855.61536
$ws.Range("N40").Value = -2394.125

$ws.Range("H93").Value = 14694.421
$ws.Range("I93").Value = 3225.5
$ws.Range("J93").Value = 34355.43
$ws.Range("K93").Value = 3225.5
$ws.Range("L93").Value = 34355.43
$ws.Range("M93").Value = -1977.5
$ws.Range("N93").Value = -36851.43

$ws.Range("H132").Value = 1603
$ws.Range("I132").Value = 1183.375
$ws.Range("J132").Value = 4000.8572
$ws.Range("K132").Value = 3550.125
$ws.Range("L132").Value = 12002.5716
$ws.Range("M132").Value = -1020.125
$ws.Range("N132").Value = -17062.5716

$ws.Range("H136").Value = 1401735
$ws.Range("I136").Value = 1839281.1
$ws.Range("K136").Value = 5517843.300000001
$ws.Range("M136").Value = -5515293.300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18456.8
$ws.Range("I81").Value = 633.3333
$ws.Range("J81").Value = 30339.111
$ws.Range("K81").Value = 1266.6666
$ws.Range("L81").Value = 60678.222
$ws.Range("M81").Value = -205.6666
$ws.Range("N81").Value = -62800.222

$ws.Range("H84").Value = 18456.8
$ws.Range("I84").Value = 633.3333
$ws.Range("J84").Value = 30339.111
$ws.Range("K84").Value = 6333.333000000001
$ws.Range("L84").Value = 303391.11
$ws.Range("M84").Value = -1029.333000000001
$ws.Range("N84").Value = -313999.11

$ws.Range("H122").Value = 1619.9166
$ws.Range("I122").Value = 1528.4286
$ws.Range("J122").Value = 1748
$ws.Range("K122").Value = 4585.2858
$ws.Range("L122").Value = 5244
$ws.Range("M122").Value = -2135.2858
$ws.Range("N122").Value = -10144

$ws.Range("H132").Value = 1048192.75
$ws.Range("I132").Value = 1439935.8
$ws.Range("J132").Value = 3544.7778
$ws.Range("K132").Value = 4319807.4
$ws.Range("L132").Value = 10634.3334
$ws.Range("M132").Value = -4317277.4
$ws.Range("N132").Value = -15694.3334

$ws.Range("H136").Value = 14606.117
$ws.Range("I136").Value = 11372.272
$ws.Range("K136").Value = 34116.81600000001
$ws.Range("M136").Value = -31566.81600000001
